$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 235-239: rotate match data (odds/results) while keeping row index (col A) fixed
# Row 235
$ws.Range("B235").Value = 6852370
$ws.Range("C235").Value = "Romania Liga I"
$ws.Range("D235").Value = "Romania Liga I"
$ws.Range("E235").Value = 45359.625
$ws.Range("F235").Value = "Dinamo Bucharest"
$ws.Range("G235").Value = "ACS UTA Batrana Doamna"
$ws.Range("H235").Value = 1
$ws.Range("I235").Value = 0
$ws.Range("J235").Value = "H"
$ws.Range("K235").Value = 2.55
$ws.Range("L235").Value = 2.875
$ws.Range("M235").Value = 3
$ws.Range("N235").Value = 2.375
$ws.Range("O235").Value = 3
$ws.Range("P235").Value = 3.1
$ws.Range("Q235").Value = -0.25
$ws.Range("R235").Value = 2
$ws.Range("S235").Value = 1.85
$ws.Range("T235").Value = 2.25
$ws.Range("U235").Value = 1.975
$ws.Range("V235").Value = 1.875
$ws.Range("W235").Value = 1.375
$ws.Range("X235").Value = -1
$ws.Range("Y235").Value = -1
$ws.Range("Z235").Value = 1
$ws.Range("AA235").Value = -1
$ws.Range("AB235").Value = -1
$ws.Range("AC235").Value = 0.875

# Row 236
$ws.Range("B236").Value = 6836277
$ws.Range("C236").Value = "Romania Liga I"
$ws.Range("D236").Value = "Romania Liga I"
$ws.Range("E236").Value = 45359.625
$ws.Range("F236").Value = "CFR Cluj"
$ws.Range("G236").Value = "AFC Hermannstadt"
$ws.Range("H236").Value = 1
$ws.Range("I236").Value = 0
$ws.Range("J236").Value = "H"
$ws.Range("K236").Value = 1.7
$ws.Range("L236").Value = 3.4
$ws.Range("M236").Value = 5
$ws.Range("N236").Value = 1.65
$ws.Range("O236").Value = 3.5
$ws.Range("P236").Value = 5.25
$ws.Range("Q236").Value = -0.75
$ws.Range("R236").Value = 1.85
$ws.Range("S236").Value = 2
$ws.Range("T236").Value = 2.25
$ws.Range("U236").Value = 1.875
$ws.Range("V236").Value = 1.975
$ws.Range("W236").Value = 0.6499999999999999
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 0.425
$ws.Range("AA236").Value = -0.5
$ws.Range("AB236").Value = -1
$ws.Range("AC236").Value = 0.9750000000000001

# Row 237
$ws.Range("B237").Value = 6870268
$ws.Range("C237").Value = "Romania Liga I"
$ws.Range("D237").Value = "Romania Liga I"
$ws.Range("E237").Value = 45359.625
$ws.Range("F237").Value = "Petrolul Ploiesti"
$ws.Range("G237").Value = "ACS Sepsi"
$ws.Range("H237").Value = 1
$ws.Range("I237").Value = 2
$ws.Range("J237").Value = "A"
$ws.Range("K237").Value = 2.8
$ws.Range("L237").Value = 3
$ws.Range("M237").Value = 2.55
$ws.Range("N237").Value = 3
$ws.Range("O237").Value = 3.2
$ws.Range("P237").Value = 2.3
$ws.Range("Q237").Value = 0.25
$ws.Range("R237").Value = 1.85
$ws.Range("S237").Value = 2
$ws.Range("T237").Value = 2.25
$ws.Range("U237").Value = 1.875
$ws.Range("V237").Value = 1.975
$ws.Range("W237").Value = -1
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = 1.3
$ws.Range("Z237").Value = -1
$ws.Range("AA237").Value = 1
$ws.Range("AB237").Value = 0.875
$ws.Range("AC237").Value = -1

# Row 238
$ws.Range("B238").Value = 6861095
$ws.Range("C238").Value = "Romania Liga I"
$ws.Range("D238").Value = "Romania Liga I"
$ws.Range("E238").Value = 45359.625
$ws.Range("F238").Value = "FC Botosani"
$ws.Range("G238").Value = "Farul Constanta"
$ws.Range("H238").Value = 0
$ws.Range("I238").Value = 0
$ws.Range("J238").Value = "D"
$ws.Range("K238").Value = 3.75
$ws.Range("L238").Value = 3.4
$ws.Range("M238").Value = 1.909
$ws.Range("N238").Value = 3.1
$ws.Range("O238").Value = 3
$ws.Range("P238").Value = 2.375
$ws.Range("Q238").Value = 0.25
$ws.Range("R238").Value = 1.775
$ws.Range("S238").Value = 2.1
$ws.Range("T238").Value = 2
$ws.Range("U238").Value = 1.8
$ws.Range("V238").Value = 2.05
$ws.Range("W238").Value = -1
$ws.Range("X238").Value = 2
$ws.Range("Y238").Value = -1
$ws.Range("Z238").Value = 0.3875
$ws.Range("AA238").Value = -0.5
$ws.Range("AB238").Value = -1
$ws.Range("AC238").Value = 1.05

# Row 239
$ws.Range("B239").Value = 6865915
$ws.Range("C239").Value = "Romania Liga I"
$ws.Range("D239").Value = "Romania Liga I"
$ws.Range("E239").Value = 45359.625
$ws.Range("F239").Value = "FC Voluntari"
$ws.Range("G239").Value = "Universitatea Cluj"
$ws.Range("H239").Value = 0
$ws.Range("I239").Value = 0
$ws.Range("J239").Value = "D"
$ws.Range("K239").Value = 3.5
$ws.Range("L239").Value = 3.25
$ws.Range("M239").Value = 2.05
$ws.Range("N239").Value = 3.4
$ws.Range("O239").Value = 3.1
$ws.Range("P239").Value = 2.15
$ws.Range("Q239").Value = 0.25
$ws.Range("R239").Value = 1.975
$ws.Range("S239").Value = 1.875
$ws.Range("T239").Value = 2.25
$ws.Range("U239").Value = 2.05
$ws.Range("V239").Value = 1.75
$ws.Range("W239").Value = -1
$ws.Range("X239").Value = 2.1
$ws.Range("Y239").Value = -1
$ws.Range("Z239").Value = 0.4875
$ws.Range("AA239").Value = -0.5
$ws.Range("AB239").Value = -1
$ws.Range("AC239").Value = 0.75

# Rows 251-256: shift data up from rows 252-257, row 257 data discarded (row deleted below)
# Row 251
$ws.Range("B251").Value = 7951750
$ws.Range("C251").Value = "Romania Liga I"
$ws.Range("D251").Value = "Romania Liga I"
$ws.Range("E251").Value = 45380.64583333334
$ws.Range("F251").Value = "ACS Sepsi"
$ws.Range("G251").Value = "CFR Cluj"
$ws.Range("K251").Value = 3.3
$ws.Range("L251").Value = 3.4
$ws.Range("M251").Value = 2.15
$ws.Range("N251").Value = 3.5
$ws.Range("O251").Value = 3.4
$ws.Range("P251").Value = 2.05
$ws.Range("Q251").Value = 0.25
$ws.Range("R251").Value = 2.05
$ws.Range("S251").Value = 1.8
$ws.Range("T251").Value = 2.25
$ws.Range("U251").Value = 1.825
$ws.Range("V251").Value = 2.025
$ws.Range("W251").Value = 0
$ws.Range("X251").Value = 0
$ws.Range("Y251").Value = 0
$ws.Range("Z251").Value = 0
$ws.Range("AA251").Value = 0

# Row 252
$ws.Range("B252").Value = 8010912
$ws.Range("C252").Value = "Romania Liga I"
$ws.Range("D252").Value = "Romania Liga I"
$ws.Range("E252").Value = 45381.44791666666
$ws.Range("F252").Value = "FC Botosani"
$ws.Range("G252").Value = "CSM Politehnica Iasi"
$ws.Range("K252").Value = 2.55
$ws.Range("L252").Value = 3.1
$ws.Range("M252").Value = 2.875
$ws.Range("N252").Value = 2.45
$ws.Range("O252").Value = 3.1
$ws.Range("P252").Value = 3
$ws.Range("Q252").Value = -0.25
$ws.Range("R252").Value = 2.1
$ws.Range("S252").Value = 1.775
$ws.Range("T252").Value = 2
$ws.Range("U252").Value = 1.8
$ws.Range("V252").Value = 2.05
$ws.Range("W252").Value = 0
$ws.Range("X252").Value = 0
$ws.Range("Y252").Value = 0
$ws.Range("Z252").Value = 0
$ws.Range("AA252").Value = 0

# Row 253
$ws.Range("B253").Value = 8010913
$ws.Range("C253").Value = "Romania Liga I"
$ws.Range("D253").Value = "Romania Liga I"
$ws.Range("E253").Value = 45381.54166666666
$ws.Range("F253").Value = "Universitatea Cluj"
$ws.Range("G253").Value = "ACS UTA Batrana Doamna"
$ws.Range("K253").Value = 1.95
$ws.Range("L253").Value = 3.4
$ws.Range("M253").Value = 4
$ws.Range("N253").Value = 1.95
$ws.Range("O253").Value = 3.4
$ws.Range("P253").Value = 4
$ws.Range("Q253").Value = -0.5
$ws.Range("R253").Value = 1.95
$ws.Range("S253").Value = 1.9
$ws.Range("T253").Value = 2.25
$ws.Range("U253").Value = 2.025
$ws.Range("V253").Value = 1.825
$ws.Range("W253").Value = 0
$ws.Range("X253").Value = 0
$ws.Range("Y253").Value = 0
$ws.Range("Z253").Value = 0
$ws.Range("AA253").Value = 0

# Row 254
$ws.Range("B254").Value = 7951749
$ws.Range("C254").Value = "Romania Liga I"
$ws.Range("D254").Value = "Romania Liga I"
$ws.Range("E254").Value = 45381.66666666666
$ws.Range("F254").Value = "CS U Craiova"
$ws.Range("G254").Value = "Rapid Bucuresti"
$ws.Range("K254").Value = 2.1
$ws.Range("L254").Value = 3.4
$ws.Range("M254").Value = 3.3
$ws.Range("N254").Value = 2.1
$ws.Range("O254").Value = 3.5
$ws.Range("P254").Value = 3.25
$ws.Range("Q254").Value = -0.25
$ws.Range("R254").Value = 1.85
$ws.Range("S254").Value = 2
$ws.Range("T254").Value = 2.5
$ws.Range("U254").Value = 1.85
$ws.Range("V254").Value = 2
$ws.Range("W254").Value = 0
$ws.Range("X254").Value = 0
$ws.Range("Y254").Value = 0
$ws.Range("Z254").Value = 0
$ws.Range("AA254").Value = 0

# Row 255
$ws.Range("B255").Value = 7951779
$ws.Range("C255").Value = "Romania Liga I"
$ws.Range("D255").Value = "Romania Liga I"
$ws.Range("E255").Value = 45382.33333333334
$ws.Range("F255").Value = "FC U Craiova 1948"
$ws.Range("G255").Value = "Otelul Galati"
$ws.Range("K255").Value = 2.3
$ws.Range("L255").Value = 3.2
$ws.Range("M255").Value = 3.2
$ws.Range("N255").Value = 2.3
$ws.Range("O255").Value = 3.2
$ws.Range("P255").Value = 3.2
$ws.Range("Q255").Value = -0.25
$ws.Range("R255").Value = 1.975
$ws.Range("S255").Value = 1.875
$ws.Range("T255").Value = 2.25
$ws.Range("U255").Value = 2.05
$ws.Range("V255").Value = 1.8
$ws.Range("W255").Value = 0
$ws.Range("X255").Value = 0
$ws.Range("Y255").Value = 0
$ws.Range("Z255").Value = 0
$ws.Range("AA255").Value = 0

# Row 256
$ws.Range("B256").Value = 7951748
$ws.Range("C256").Value = "Romania Liga I"
$ws.Range("D256").Value = "Romania Liga I"
$ws.Range("E256").Value = 45382.625
$ws.Range("F256").Value = "Farul Constanta"
$ws.Range("G256").Value = "FCSB"
$ws.Range("K256").Value = 3.6
$ws.Range("L256").Value = 3.3
$ws.Range("M256").Value = 2
$ws.Range("N256").Value = 3.6
$ws.Range("O256").Value = 3.3
$ws.Range("P256").Value = 2.05
$ws.Range("Q256").Value = 0.25
$ws.Range("R256").Value = 2.05
$ws.Range("S256").Value = 1.8
$ws.Range("T256").Value = 2.25
$ws.Range("U256").Value = 1.825
$ws.Range("V256").Value = 2.025
$ws.Range("W256").Value = 0
$ws.Range("X256").Value = 0
$ws.Range("Y256").Value = 0
$ws.Range("Z256").Value = 0
$ws.Range("AA256").Value = 0

# Remove the now-obsolete last row (257), which has been fully absorbed into rows above
$ws.Rows.Item(257).Delete()
